# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh updates to the Jenova_Profits workbook
# (columns H/I/J/K/L/M/N = currentAveragePrice*, LevePrice*, LeveProfit* per leve row)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2053.9836
$ws.Range("J17").Value = 2053.9836
$ws.Range("L17").Value = 6161.950800000001
$ws.Range("N17").Value = -6497.950800000001
$ws.Range("H28").Value = 56706.61
$ws.Range("I28").Value = 56706.61
$ws.Range("K28").Value = 56706.61
$ws.Range("M28").Value = -56221.61
$ws.Range("H33").Value = 852.93335
$ws.Range("I33").Value = 852.93335
$ws.Range("K33").Value = 852.93335
$ws.Range("M33").Value = -623.93335
$ws.Range("H40").Value = 8772.637000000001
$ws.Range("J40").Value = 10749.167
$ws.Range("L40").Value = 10749.167
$ws.Range("N40").Value = -11099.167
$ws.Range("H62").Value = 13895110
$ws.Range("I62").Value = 41669864
$ws.Range("K62").Value = 41669864
$ws.Range("M62").Value = -41669240
$ws.Range("H65").Value = 13895110
$ws.Range("I65").Value = 41669864
$ws.Range("K65").Value = 208349320
$ws.Range("M65").Value = -208346200
$ws.Range("H111").Value = 63092.062
$ws.Range("I111").Value = 83844.164
$ws.Range("J111").Value = 835.75
$ws.Range("K111").Value = 251532.492
$ws.Range("L111").Value = 2507.25
$ws.Range("M111").Value = -248465.492
$ws.Range("N111").Value = -8641.25
$ws.Range("H129").Value = 1919.4
$ws.Range("I129").Value = 545.875
$ws.Range("J129").Value = 2835.0833
$ws.Range("K129").Value = 1637.625
$ws.Range("L129").Value = 8505.249899999999
$ws.Range("M129").Value = 3362.375
$ws.Range("N129").Value = -18505.2499
$ws.Range("H132").Value = 13189.194
$ws.Range("I132").Value = 2087.8438
$ws.Range("J132").Value = 102000
$ws.Range("K132").Value = 6263.5314
$ws.Range("L132").Value = 306000
$ws.Range("M132").Value = -3733.5314
$ws.Range("N132").Value = -311060
$ws.Range("H133").Value = 42241.56
$ws.Range("J133").Value = 42241.56
$ws.Range("L133").Value = 42241.56
$ws.Range("N133").Value = -52361.56
$ws.Range("H135").Value = 716425.0600000001
$ws.Range("I135").Value = 771073.1
$ws.Range("K135").Value = 6939657.899999999
$ws.Range("M135").Value = -6937122.899999999
$ws.Range("H137").Value = 3937.9697
$ws.Range("I137").Value = 4124.3335
$ws.Range("K137").Value = 12373.0005
$ws.Range("M137").Value = -9823.000499999998
$ws.Range("H138").Value = 3635.9048
$ws.Range("I138").Value = 1557.2632
$ws.Range("K138").Value = 4671.7896
$ws.Range("M138").Value = 468.2103999999999
$ws.Range("H141").Value = 2847.3333
$ws.Range("I141").Value = 2511.55
$ws.Range("J141").Value = 4526.25
$ws.Range("K141").Value = 7534.650000000001
$ws.Range("L141").Value = 13578.75
$ws.Range("M141").Value = -2354.650000000001
$ws.Range("N141").Value = -23938.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3912.7014
$ws.Range("I32").Value = 3557.7812
$ws.Range("K32").Value = 3557.7812
$ws.Range("M32").Value = -3270.7812
$ws.Range("H45").Value = 2621.4375
$ws.Range("I45").Value = 1640.4546
$ws.Range("K45").Value = 1640.4546
$ws.Range("M45").Value = -1263.4546
$ws.Range("H74").Value = 1647.0555
$ws.Range("I74").Value = 1581.9286
$ws.Range("J74").Value = 1875
$ws.Range("K74").Value = 1581.9286
$ws.Range("L74").Value = 1875
$ws.Range("M74").Value = -707.9286
$ws.Range("N74").Value = -3623
$ws.Range("H77").Value = 1647.0555
$ws.Range("I77").Value = 1581.9286
$ws.Range("J77").Value = 1875
$ws.Range("K77").Value = 7909.643
$ws.Range("L77").Value = 9375
$ws.Range("M77").Value = -3541.643
$ws.Range("N77").Value = -18111
$ws.Range("H132").Value = 3668.6863
$ws.Range("I132").Value = 3425.0466
$ws.Range("J132").Value = 4978.25
$ws.Range("K132").Value = 10275.1398
$ws.Range("L132").Value = 14934.75
$ws.Range("M132").Value = -7745.139800000001
$ws.Range("N132").Value = -19994.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = -127
$ws.Range("H94").Value = 2241.6428
$ws.Range("J94").Value = 486
$ws.Range("L94").Value = 486
$ws.Range("N94").Value = -1388
$ws.Range("H132").Value = 50000
$ws.Range("J132").Value = 50000
$ws.Range("L132").Value = 50000
$ws.Range("N132").Value = -60120
$ws.Range("H134").Value = 30857.361
$ws.Range("I134").Value = 2966.647
$ws.Range("K134").Value = 8899.940999999999
$ws.Range("M134").Value = -6364.940999999999
$ws.Range("H137").Value = 49999.5
$ws.Range("J137").Value = 49999.5
$ws.Range("L137").Value = 49999.5
$ws.Range("N137").Value = -60199.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 39189.57
$ws.Range("I31").Value = 1071.6666
$ws.Range("J31").Value = 67778
$ws.Range("K31").Value = 1071.6666
$ws.Range("L31").Value = 67778
$ws.Range("M31").Value = -776.6666
$ws.Range("N31").Value = -68368
$ws.Range("H34").Value = 39189.57
$ws.Range("I34").Value = 1071.6666
$ws.Range("J34").Value = 67778
$ws.Range("K34").Value = 1071.6666
$ws.Range("L34").Value = 67778
$ws.Range("M34").Value = -869.6666
$ws.Range("N34").Value = -68182
$ws.Range("H105").Value = 1008.7143
$ws.Range("J105").Value = 1007
$ws.Range("L105").Value = 1007
$ws.Range("N105").Value = -4501
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H132").Value = 2634.75
$ws.Range("I132").Value = 2500
$ws.Range("K132").Value = 7500
$ws.Range("M132").Value = -4970
$ws.Range("H134").Value = 502717.44
$ws.Range("I134").Value = 2860.4736
$ws.Range("K134").Value = 8581.4208
$ws.Range("M134").Value = -6046.4208

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 182.29411
$ws.Range("I12").Value = 80
$ws.Range("J12").Value = 204.21428
$ws.Range("K12").Value = 240
$ws.Range("L12").Value = 612.64284
$ws.Range("M12").Value = -67
$ws.Range("N12").Value = -958.64284
$ws.Range("H86").Value = 1178.4
$ws.Range("I86").Value = 296
$ws.Range("J86").Value = 1766.6666
$ws.Range("K86").Value = 888
$ws.Range("L86").Value = 5299.9998
$ws.Range("M86").Value = 298
$ws.Range("N86").Value = -7671.9998
$ws.Range("H89").Value = 1178.4
$ws.Range("I89").Value = 296
$ws.Range("J89").Value = 1766.6666
$ws.Range("K89").Value = 2664
$ws.Range("L89").Value = 15899.9994
$ws.Range("M89").Value = 3264
$ws.Range("N89").Value = -27755.9994
$ws.Range("H133").Value = 7470.75
$ws.Range("I133").Value = 7036.1665
$ws.Range("K133").Value = 21108.4995
$ws.Range("M133").Value = -16048.4995
$ws.Range("H138").Value = 3565.4443
$ws.Range("I138").Value = 3298.8333
$ws.Range("J138").Value = 4098.6665
$ws.Range("K138").Value = 9896.499899999999
$ws.Range("L138").Value = 12295.9995
$ws.Range("M138").Value = -4756.499899999999
$ws.Range("N138").Value = -22575.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 958.8946999999999
$ws.Range("I97").Value = 961.4375
$ws.Range("J97").Value = 945.3333
$ws.Range("K97").Value = 961.4375
$ws.Range("L97").Value = 945.3333
$ws.Range("M97").Value = -465.4375
$ws.Range("N97").Value = -1937.3333
$ws.Range("H126").Value = 3417.0454
$ws.Range("I126").Value = 3061.5
$ws.Range("K126").Value = 9184.5
$ws.Range("M126").Value = -6714.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 756.2857
$ws.Range("J22").Value = 831.6667
$ws.Range("L22").Value = 831.6667
$ws.Range("N22").Value = -1421.6667
$ws.Range("H27").Value = 756.2857
$ws.Range("J27").Value = 831.6667
$ws.Range("L27").Value = 831.6667
$ws.Range("N27").Value = -1045.6667
$ws.Range("H46").Value = 2379.7
$ws.Range("I46").Value = 2549.6667
$ws.Range("K46").Value = 2549.6667
$ws.Range("M46").Value = -2361.6667
$ws.Range("H127").Value = 94489
$ws.Range("J127").Value = 94489
$ws.Range("L127").Value = 94489
$ws.Range("N127").Value = -104409
$ws.Range("H136").Value = 838157.8
$ws.Range("I136").Value = 838157.8
$ws.Range("K136").Value = 2514473.4
$ws.Range("M136").Value = -2511923.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 719.4815
$ws.Range("I107").Value = 762.8182
$ws.Range("K107").Value = 2288.4546
$ws.Range("M107").Value = -368.4546
$ws.Range("H122").Value = 27029010
$ws.Range("I122").Value = 31251660
$ws.Range("J122").Value = 4060
$ws.Range("K122").Value = 93754980
$ws.Range("L122").Value = 12180
$ws.Range("M122").Value = -93752530
$ws.Range("N122").Value = -17080
$ws.Range("H126").Value = 1597
$ws.Range("I126").Value = 1204
$ws.Range("K126").Value = 3612
$ws.Range("M126").Value = -1142
$ws.Range("H132").Value = 96257.55
$ws.Range("I132").Value = 3888.8333
$ws.Range("K132").Value = 11666.4999
$ws.Range("M132").Value = -9136.499899999999
$ws.Range("H136").Value = 12423000
$ws.Range("I136").Value = 15154864
$ws.Range("K136").Value = 45464592
$ws.Range("M136").Value = -45462042
